{"js": "// Replace each two-digit-divided-by-one-digit division expression in the\n// worksheet table with its updated value, per the commit's regenerated\n// problem set. Every source expression string is unique in the document,\n// so a direct search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"81\u00f77=\", \"72\u00f79=\"],\n  [\"59\u00f75=\", \"35\u00f79=\"],\n  [\"43\u00f75=\", \"88\u00f79=\"],\n  [\"34\u00f72=\", \"35\u00f74=\"],\n  [\"77\u00f79=\", \"80\u00f75=\"],\n  [\"58\u00f76=\", \"20\u00f74=\"],\n  [\"63\u00f79=\", \"74\u00f76=\"],\n  [\"99\u00f76=\", \"55\u00f72=\"],\n  [\"95\u00f74=\", \"95\u00f75=\"],\n  [\"72\u00f77=\", \"47\u00f79=\"],\n  [\"62\u00f76=\", \"16\u00f73=\"],\n  [\"18\u00f73=\", \"34\u00f74=\"],\n  [\"98\u00f75=\", \"51\u00f78=\"],\n  [\"27\u00f72=\", \"99\u00f75=\"],\n  [\"43\u00f79=\", \"17\u00f76=\"],\n  [\"51\u00f75=\", \"36\u00f77=\"],\n  [\"92\u00f74=\", \"37\u00f78=\"],\n  [\"54\u00f74=\", \"88\u00f75=\"],\n  [\"11\u00f73=\", \"13\u00f72=\"],\n  [\"82\u00f72=\", \"58\u00f72=\"],\n  [\"59\u00f73=\", \"54\u00f78=\"],\n  [\"22\u00f73=\", \"76\u00f79=\"],\n  [\"55\u00f78=\", \"49\u00f73=\"],\n  [\"85\u00f78=\", \"87\u00f74=\"],\n  [\"71\u00f76=\", \"55\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-divided-by-one-digit division expression in the\n# worksheet table with its updated value, per the commit's regenerated\n# problem set. Every source expression string is unique in the document,\n# so Find/Execute with wdReplaceAll (2) on each pair is unambiguous.\n$pairs = @(\n  @(\"81\u00f77=\", \"72\u00f79=\"),\n  @(\"59\u00f75=\", \"35\u00f79=\"),\n  @(\"43\u00f75=\", \"88\u00f79=\"),\n  @(\"34\u00f72=\", \"35\u00f74=\"),\n  @(\"77\u00f79=\", \"80\u00f75=\"),\n  @(\"58\u00f76=\", \"20\u00f74=\"),\n  @(\"63\u00f79=\", \"74\u00f76=\"),\n  @(\"99\u00f76=\", \"55\u00f72=\"),\n  @(\"95\u00f74=\", \"95\u00f75=\"),\n  @(\"72\u00f77=\", \"47\u00f79=\"),\n  @(\"62\u00f76=\", \"16\u00f73=\"),\n  @(\"18\u00f73=\", \"34\u00f74=\"),\n  @(\"98\u00f75=\", \"51\u00f78=\"),\n  @(\"27\u00f72=\", \"99\u00f75=\"),\n  @(\"43\u00f79=\", \"17\u00f76=\"),\n  @(\"51\u00f75=\", \"36\u00f77=\"),\n  @(\"92\u00f74=\", \"37\u00f78=\"),\n  @(\"54\u00f74=\", \"88\u00f75=\"),\n  @(\"11\u00f73=\", \"13\u00f72=\"),\n  @(\"82\u00f72=\", \"58\u00f72=\"),\n  @(\"59\u00f73=\", \"54\u00f78=\"),\n  @(\"22\u00f73=\", \"76\u00f79=\"),\n  @(\"55\u00f78=\", \"49\u00f73=\"),\n  @(\"85\u00f78=\", \"87\u00f74=\"),\n  @(\"71\u00f76=\", \"55\u00f73=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
